# Fruta / hortaliza, semanal
# Insert a new weekly pair of rows (Primera/Segunda) at the top of the
# Brocoli data block (row 369) for Terminal La Palmera de La Serena.
# Everything below shifts down by two rows; the old last pair (formerly
# rows 419/420) lands on the newly-grown rows 421/422.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 369:420 down by two rows, duplicating row 369's formatting
# into the freshly inserted rows.
$ws.Rows("369:370").Insert()

# Common (unchanged) column values shared by every data row in this block.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$catId     = 100112023
$categoria = "Brócoli"
$variedad  = "Sin especificar"
$unidad    = "`$/unidad"
$origen    = "Provincia del Elquí"
$kgUnidad  = 1
$clasif    = "Hortaliza"

# New "Primera" record (row 369).
$ws.Cells.Item(369, 1).Value  = $mercadoId
$ws.Cells.Item(369, 2).Value  = $mercado
$ws.Cells.Item(369, 3).Value  = $region
$ws.Cells.Item(369, 4).Value  = 44491
$ws.Cells.Item(369, 5).Value  = $codreg
$ws.Cells.Item(369, 6).Value  = $catId
$ws.Cells.Item(369, 7).Value  = $categoria
$ws.Cells.Item(369, 8).Value  = $variedad
$ws.Cells.Item(369, 9).Value  = "Primera"
$ws.Cells.Item(369, 10).Value = 3400
$ws.Cells.Item(369, 11).Value = 600
$ws.Cells.Item(369, 12).Value = 700
$ws.Cells.Item(369, 13).Value = 650
$ws.Cells.Item(369, 14).Value = $unidad
$ws.Cells.Item(369, 15).Value = $origen
$ws.Cells.Item(369, 16).Value = 650
$ws.Cells.Item(369, 17).Value = $kgUnidad
$ws.Cells.Item(369, 18).Value = $clasif

# New "Segunda" record (row 370).
$ws.Cells.Item(370, 1).Value  = $mercadoId
$ws.Cells.Item(370, 2).Value  = $mercado
$ws.Cells.Item(370, 3).Value  = $region
$ws.Cells.Item(370, 4).Value  = 44491
$ws.Cells.Item(370, 5).Value  = $codreg
$ws.Cells.Item(370, 6).Value  = $catId
$ws.Cells.Item(370, 7).Value  = $categoria
$ws.Cells.Item(370, 8).Value  = $variedad
$ws.Cells.Item(370, 9).Value  = "Segunda"
$ws.Cells.Item(370, 10).Value = 1600
$ws.Cells.Item(370, 11).Value = 500
$ws.Cells.Item(370, 12).Value = 550
$ws.Cells.Item(370, 13).Value = 525
$ws.Cells.Item(370, 14).Value = $unidad
$ws.Cells.Item(370, 15).Value = $origen
$ws.Cells.Item(370, 16).Value = 525
$ws.Cells.Item(370, 17).Value = $kgUnidad
$ws.Cells.Item(370, 18).Value = $clasif
